$wb = $excel.ActiveWorkbook

# Move "PushMsgService(消息推送中心)" so it sits right before "appBackEnd(手机app后台)",
# i.e. swap the order of the 2nd and 3rd tabs.
$pushMsg = $wb.Worksheets.Item("PushMsgService(消息推送中心)")
$appBackEnd = $wb.Worksheets.Item("appBackEnd(手机app后台)")
$pushMsg.Move($appBackEnd) | Out-Null

# Rename the three sheets to their new product names.
$wb.Worksheets.Item("idm(用户信息管理系统)").Name = "ACV-UA"
$wb.Worksheets.Item("PushMsgService(消息推送中心)").Name = "ACV-NA"
$wb.Worksheets.Item("appBackEnd(手机app后台)").Name = "ACV-VA"

# Restore the per-sheet cursor positions noted in the saved file.
$wb.Worksheets.Item("ACV-UA").Activate() | Out-Null
$wb.Worksheets.Item("ACV-UA").Range("B23").Select() | Out-Null

$wb.Worksheets.Item("ACV-NA").Activate() | Out-Null
$wb.Worksheets.Item("ACV-NA").Range("C23").Select() | Out-Null

$wb.Worksheets.Item("ACV-VA").Activate() | Out-Null
$wb.Worksheets.Item("ACV-VA").Range("C19").Select() | Out-Null
